# Add new "Kredit / Vorhaben" vocabulary entries to the "vocab" sheet.
# Mirrors 14 new rows (169-182) appended after the existing data (which
# ended at row 168), using the same column layout:
#   A = German, B = English, C = lesson date, D = lesson number, E = phrase/word

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vocab")

$lastRow = 168
$lessonDate = "3/16/2022"   # serial 44636
$lessonNumber = 7

$rows = @(
    @("einen Kredit aufnehmen", "to take out a loan", "phrase"),
    @("einen Kredit genehmigen", "to grant a loan", "phrase"),
    @("einen Kredit gewähren", "to grant a loan", "phrase"),
    @("einen Kredit (bekommen) erhalten", "to get a loan", "phrase"),
    @("einen Kredit verweigern", "to reject a credit", "phrase"),
    @("die Kreditwürdigkeit / die finanzielle Lage prüfen [überprüfen] / einschätzen / bewerten", "to check / assess the financial situation", "phrase"),
    @("jemand vom Erfolg eines Vorhabens überzeugen", "to convince someone of the success of a project", "phrase"),
    @("einen Businessplan / Finanzierungsplan erstellen / machen / erklären / darlegen", "to prepare / explain / present a business plan", "phrase"),
    @("in ein Vorhaben investieren", "to invest in a project", "phrase"),
    @("Fristen setzen / einhalten", "to set a timeline / to meet a timeline", "word"),
    @("der Gegenstand / die Gegenstände", "object / items", "word"),
    @("das Vermögen", "the fortune / assets", "word"),
    @("die Gelegenheit", "opportunity", "word"),
    @("der Nutzen", "the need / the result / profit / dividend?", "word")
)

$r = $lastRow
foreach ($row in $rows) {
    $r = $r + 1

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]

    # Copy the date cell's format from the row above so the new cell gets
    # the same short-date number format/style as the rest of column C.
    $ws.Cells.Item($r - 1, 3).Copy()
    $ws.Cells.Item($r, 3).PasteSpecial(-4122)
    $ws.Cells.Item($r, 3).Value = $lessonDate

    $ws.Cells.Item($r, 4).Value = $lessonNumber
    $ws.Cells.Item($r, 5).Value = $row[2]
}

$excel.CutCopyMode = 0

# Update the window scroll position / selection to match the saved state.
$ws.Range("A150").Select()
$ws.Range("E177").Select()
